# 9th Stab - Cosmetic Changes
# Insert two new "week" columns (Jun_17, Jun_15) in front of the existing
# Jun_13 / Jun_10 columns, pushing the existing data two columns to the
# right (B->D, C->E) and seeding the two new columns with the same "UN"
# placeholder used for unchanged analysts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns B (Jun_13) and C (Jun_10) right by two, inserting
# two blank columns at B:C.
$ws.Columns("B:C").Insert()

# New header row values - set C1 before B1 so the shared-string table picks
# up "Jun_15" ahead of "Jun_17", matching the original append order.
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Fill the two new week columns with the "UN" (unchanged) placeholder used
# throughout the rest of the sheet.
$ws.Range("B2:C27").Value = "UN"

# Keep the newly inserted columns the same width as the original data
# column (8 characters).
$ws.Columns("C").ColumnWidth = 7.166666666666667
$ws.Columns("D").ColumnWidth = 7.166666666666667
